$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-format cells whose new value would otherwise be auto-detected as a
# number (losing e.g. trailing zeros) so they round-trip as plain text,
# matching the source data which is always stored as text.
$textCells = @('D5', 'D8', 'D9', 'D11', 'D18', 'D20', 'D22', 'D24', 'D25', 'D26', 'D31', 'D32', 'D36', 'D37', 'D39', 'D44', 'D48', 'D51', 'D41', 'D42')
foreach ($c in $textCells) {
    $ws.Range($c).NumberFormat = "@"
}

$ws.Range('D2').Value = '33.642.96'
$ws.Range('E2').Value = '  -0.95%  '
$ws.Range('D3').Value = '1.769.13'
$ws.Range('E3').Value = '  -0.89%  '
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').Value = '223.44'
$ws.Range('E5').Value = '  +0.79%  '
$ws.Range('E6').Value = '  -1.13%  '
$ws.Range('E7').Value = '  +0.05%  '
$ws.Range('D8').Value = '31.76'
$ws.Range('E8').Value = '  +1.27%  '
$ws.Range('D9').Value = '0.289'
$ws.Range('E9').Value = '  +1.72%  '
$ws.Range('E10').Value = '  -3.82%  '
$ws.Range('D11').Value = '0.0934'
$ws.Range('E11').Value = '  +1.31%  '
$ws.Range('D12').Value = '2.024.27'
$ws.Range('E12').Value = '  -0.81%  '
$ws.Range('E13').Value = '  +4.18%  '
$ws.Range('D14').Value = '1.772.95'
$ws.Range('E14').Value = '  -0.69%  '
$ws.Range('D15').Value = '33.685.06'
$ws.Range('E15').Value = '  -0.78%  '
$ws.Range('E16').Value = '  -3.24%  '
$ws.Range('D18').Value = '66.30'
$ws.Range('E18').Value = '  -2.46%  '
$ws.Range('D19').Value = '0.0₃0773'
$ws.Range('E19').Value = '  -1.17%  '
$ws.Range('D20').Value = '237.54'
$ws.Range('E20').Value = '  -2.99%  '
$ws.Range('E21').Value = '  +0.03%  '
$ws.Range('D22').Value = '10.52'
$ws.Range('E22').Value = '  -1.65%  '
$ws.Range('E23').Value = '  -1.92%  '
$ws.Range('D24').Value = '2.05'
$ws.Range('E24').Value = '  -2.27%  '
$ws.Range('D25').Value = '159.30'
$ws.Range('E25').Value = '  +1.06%  '
$ws.Range('D26').Value = '16.05'
$ws.Range('E26').Value = '  -1.94%  '
$ws.Range('E27').Value = '  -0.02%  '
$ws.Range('E28').Value = '  -0.45%  '
$ws.Range('E29').Value = '  +0.12%  '
$ws.Range('E30').Value = '  +1.20%  '
$ws.Range('D31').Value = '0.0510'
$ws.Range('E31').Value = '  -1.68%  '
$ws.Range('D32').Value = '3.58'
$ws.Range('E32').Value = '  -2.87%  '
$ws.Range('E33').Value = '  -0.48%  '
$ws.Range('E34').Value = '  -1.63%  '
$ws.Range('D35').Value = '1.377.55'
$ws.Range('E35').Value = '  -2.33%  '
$ws.Range('D36').Value = '0.644'
$ws.Range('E36').Value = '  +0.46%  '
$ws.Range('D37').Value = '1.02'
$ws.Range('E37').Value = '  -2.54%  '
$ws.Range('E38').Value = '  -1.26%  '
$ws.Range('D39').Value = '2.22'
$ws.Range('E39').Value = '  +5.36%  '
$ws.Range('E40').Value = '  +0.84%  '
$ws.Range('E43').Value = '  -1.98%  '
$ws.Range('D44').Value = '13.50'
$ws.Range('E44').Value = '  +14.44%  '
$ws.Range('E45').Value = '  +4.06%  '
$ws.Range('D46').Value = '0.0₆0137'
$ws.Range('E46').Value = '  +13.98%  '
$ws.Range('E47').Value = '  +1.09%  '
$ws.Range('D48').Value = '106.75'
$ws.Range('E48').Value = '  +1.14%  '
$ws.Range('E49').Value = '  -2.23%  '
$ws.Range('D50').Value = '1.924.61'
$ws.Range('E50').Value = '  -0.46%  '
$ws.Range('D51').Value = '0.999'
$ws.Range('E51').Value = '  +0.20%  '
$ws.Range('B41').Value = 'Aave'
$ws.Range('C41').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D41').Value = '77.62'
$ws.Range('E41').Value = '  -2.47%  '
$ws.Range('B42').Value = 'ARBITRUM'
$ws.Range('C42').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D42').Value = '0.902'
$ws.Range('E42').Value = '  -3.68%  '
